$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows at the top; this shifts the existing header/data
#    rows (1,2) down to (3,4), carrying their hyperlink reference intact.
$ws.Rows("1:2").Insert()

# 2) Rewrite the (now shifted) header row with the new lower-case column
#    names.
$ws.Range("A3").Value = "nome"
$ws.Range("B3").Value = "edv"
$ws.Range("C3").Value = "trilha"
$ws.Range("D3").Value = "gestor"
$ws.Range("E3").Value = "gestor_email"

# 3) Rewrite the (now shifted) example-employee row.
$ws.Range("A4").Value = "Lucas"
$ws.Range("B4").Value = 12345678
$ws.Range("C4").Value = "Lider"
$ws.Range("D4").Value = "Henrique Dona"
$ws.Range("E4").Value = "dona@br.bosch.com"

# 4) Add a new second example row.
$ws.Range("A5").Value = "Giovana"
$ws.Range("B5").Value = 87654321
$ws.Range("C5").Value = "Qualidade"
$ws.Range("D5").Value = "Henrique Dona"
$ws.Range("E5").Value = "dona@br.bosch.com"

# 5) Put the big instructional banner text in A1 (merged A1:E2 later).
$ws.Range("A1").Value = "ESSA TABELA SERVE APENAS DE EXEMPLO, POR FAVOR SIGA O PADRÃO DELA `n(AS TRILHAS DEVEM TER A PRIMEIRA LETRA EM MAIÚSCULO, EDV DEVE TER EXATOS 8 NÚMEROS ASSIM COMO NO CRACHA)"
$ws.Rows("1:1").AutoFit()

# 6) Re-point the original hyperlink (it stayed on "E2" after the row
#    insert) to its new cell, then add the new hyperlink for row 5.
$ws.Hyperlinks(1).Range = $ws.Range("E4")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:dona@br.bosch.com")

# 7) Styling: row 5's name/edv/trilha cells get horizontal-center-only
#    alignment (built on a scratch cell first, then stamped across with a
#    single paste so only one new style entry is produced).
$scratch = $ws.Range("H100")
$scratch.HorizontalAlignment = -4108
$scratch.Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$scratch.ClearContents()
$excel.CutCopyMode = $false

# 8) Styling: the red banner area (B1:E1, part of the merged A1:E2 block)
#    gets a solid red fill with white centered text; A1 gets the same
#    plus word-wrap. Build each on its own scratch cell first so exactly
#    one new style gets created per look, then paste across in one shot.
$scratch2 = $ws.Range("H101")
$scratch2.Interior.Color = 255
$scratch2.Font.ThemeColor = 2
$scratch2.HorizontalAlignment = -4108
$scratch2.VerticalAlignment = -4108
$scratch2.Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$scratch2.ClearContents()
$excel.CutCopyMode = $false

$scratch3 = $ws.Range("H102")
$scratch3.Interior.Color = 255
$scratch3.Font.ThemeColor = 2
$scratch3.HorizontalAlignment = -4108
$scratch3.VerticalAlignment = -4108
$scratch3.WrapText = $true
$scratch3.Copy()
$ws.Range("A1").PasteSpecial(-4122)
$scratch3.ClearContents()
$excel.CutCopyMode = $false

# 9) Merge the banner cells and select it (mirrors the authored sheet).
$ws.Range("A1:E2").Merge()
$ws.Range("A1:E2").Select()

# 10) Page setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
